# Add the missing 2017 year of data to the "BSoAIGtAP" sheet.
#
# Real-world edit: the author inserted a new column before column B on the
# "BSoAIGtAP" worksheet (shifting 2018..2050 one column to the right, to
# C..AI) and filled the new column B with the 2017 year label/value:
#   B1 = 2017
#   B2 = formula mirroring the (now shifted) "Data!A33" cell, i.e. =$C2
# The sheet then became the active sheet/tab with cell L13 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BSoAIGtAP")

# Insert a new blank column at B; everything in B:AH shifts right to C:AI
# (formulas referencing the old column B are auto-adjusted by Excel to
# point at the new location, column C).
$rng = $ws.Range("B1:B2")
$rng.EntireColumn.Insert()

# Fill in the newly-inserted 2017 column.
$ws.Cells.Item(1, 2).Value = 2017
$ws.Cells.Item(2, 2).Formula = "=`$C2"

# Make this the active sheet, with the same cell selected as in the saved
# file (L13).
$ws.Activate()
$ws.Range("L13").Select()
